$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto price/volume snapshot values (D: Price, E: Volume(1h)).
# Values that look numeric must be forced to Text (leading-apostrophe,
# like typing into Excel's formula bar) so they keep their original
# display formatting, then the cell style is reset to Normal so no
# extraneous "quote prefix" number-format style is introduced.

$ws.Range("D2").Value = "29.512.60"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.978.13"
$ws.Range("E3").Value = "  +4.01%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").Value = "'327.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "'0.4662"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").Value = "'0.07957"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").Value = "'0.9938"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").Value = "'22.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.90%  "
$ws.Range("D12").Value = "1.978.31"
$ws.Range("E12").Value = "  +4.04%  "
$ws.Range("D13").Value = "'7.203"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("D14").Value = "'5.852"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("D15").Value = "'0.07091"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").Value = "'87.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "'0.000009948"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "29.513.25"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("D22").Value = "'5.562"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.53%  "
$ws.Range("D23").Value = "'11.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("D24").Value = "2.223.13"
$ws.Range("E24").Value = "  +4.23%  "
$ws.Range("D25").Value = "'2.110"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").Value = "'158.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("D27").Value = "'19.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("D28").Value = "'5.790"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.22%  "
$ws.Range("D29").Value = "'119.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").Value = "'0.09425"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("D32").Value = "'0.8986"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("D33").Value = "'5.239"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "'1.324"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").Value = "'3.193"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("D36").Value = "'0.05821"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("D37").Value = "'1.172"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.41%  "
$ws.Range("D38").Value = "'0.02104"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("D39").Value = "'7.794"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").Value = "'0.5723"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").Value = "'0.000003140"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +44.75%  "
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("D44").Value = "'2.797"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.74%  "
$ws.Range("D45").Value = "'11.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").Value = "'0.5371"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").Value = "'2.185"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").Value = "'0.06930"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("D49").Value = "'114.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("D50").Value = "'1.831"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "'0.3041"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.69%  "
